# AverageTimes.xlsx edit: add a Java/Python/C++ breakdown next to the existing
# sort-algorithm timings, mirror it as a second (pivoted) table further down
# the sheet, and extend the stacked bar chart with two more (currently empty)
# series so it's ready to be filled in with the Python / C++ numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New label column (E) next to the existing Selection/Bubble/Merge/Quick
#     Sort timings, marking which language each row of timings belongs to.
$ws.Range("E2").Value = "Java"
$ws.Range("E3").Value = "Python"
$ws.Range("E4").Value = "C++"

# --- A second, pivoted copy of the table starting at row 21: algorithm
#     names across the header row, language names down column A.
$ws.Range("B21").Value = "Selection Sort"
$ws.Range("C21").Value = "Bubble Sort"
$ws.Range("D21").Value = "Merge Sort"
$ws.Range("E21").Value = "Quick Sort"

$ws.Range("A22").Value = "Java"
$ws.Range("B22").Value = 5.5
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 14.7
$ws.Range("E22").Value = 11.9

$ws.Range("A23").Value = "Python"
$ws.Range("A24").Value = "C++"

# --- Column widths: B-D grow a bit (no longer auto bestFit) and the new
#     column E gets an explicit custom width too.
$ws.Columns.Item(2).ColumnWidth = 12.333333333333332
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 9.666666666666666
$ws.Columns.Item(5).ColumnWidth = 9.333333333333332

# --- Chart: the existing stacked bar chart only had one series (row 2 -
#     the Java timings). Add two more series, stacked on top, pointing at
#     rows 3 and 4 (Python / C++ - still empty, to be filled in later).
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()

$s2 = $sc.NewSeries()
$s2.Formula = "=SERIES(,Sheet1!`$A`$1:`$D`$1,Sheet1!`$A`$3:`$D`$3,2)"

$s3 = $sc.NewSeries()
$s3.Formula = "=SERIES(,Sheet1!`$A`$1:`$D`$1,Sheet1!`$A`$4:`$D`$4,3)"

# --- Leave the selection where the author was last working.
$null = $ws.Range("D16").Select()
